# SCD_JobPackage.xlsx edit script
# - Renumbers Question_Concept_Code values (column E) for rows 137-148
# - Sets the Category (column G) for row 149 and updates its Question_Concept_Code
# - Appends 12 new rows (150-161) with new concepts/categories

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix up existing rows 137-148 (column E renumbering) ---
$eFixes = @{
    137 = 2000000154
    138 = 2000000155
    139 = 2000000156
    140 = 2000000157
    141 = 2000000158
    142 = 2000000159
    143 = 2000000160
    144 = 2000000161
    145 = 2000000162
    146 = 2000000163
    147 = 2000000164
    148 = 2000000165
}

foreach ($r in $eFixes.Keys) {
    $ws.Cells.Item($r, 5).Value = $eFixes[$r]
}

# --- Row 149: update Question_Concept_Code and set Category ---
$ws.Cells.Item(149, 5).Value = 2000000166
$ws.Cells.Item(149, 6).Value = "Treatment"
$ws.Cells.Item(149, 7).Value = "Primary Treatment"

# --- New rows 150-161 ---
$newRows = @(
    @(2000000167, "First Sickle Cell Encounter", "SCD Registry", "Observation", 2000000167, "Diagnoses", "Primary Diagnoses"),
    @(2000000168, "Last Sickle Cell Encounter", "SCD Registry", "Observation", 2000000168, "Diagnoses", "Primary Diagnoses"),
    @(2000000169, "Adakveo", "SCD Registry", "Observation", 2000000169, "SCD Medication", "Medication"),
    @(2000000170, "Oxbrtya", "SCD Registry", "Observation", 2000000170, "SCD Medication", "Medication"),
    @(2000000171, "Voxeletor", "SCD Registry", "Observation", 2000000171, "SCD Medication", "Medication"),
    @(2000000172, "Hydroxyurea", "SCD Registry", "Observation", 2000000172, "SCD Medication", "Medication"),
    @(2000000173, "Endari", "SCD Registry", "Observation", 2000000173, "SCD Medication", "Medication"),
    @(2000000174, "Deferasirox", "SCD Registry", "Observation", 2000000174, "SCD Medication", "Medication"),
    @(2000000175, "Date of most recent transfusion of blood product", "SCD Registry", "Observation", 2000000175, "Other Medical History", "Transfusion"),
    @(2000000176, "Medication", "SCD Registry", "Observation", 2000000176, "Other Medical History", "Medication"),
    @(2000000177, "Primary Care Provider Name", "Common Registry", "Observation", 2000000177, "Demographics", ""),
    @(2000000178, "Sickle Cell Diagnosis", "SCD Registry", "Observation", 2000000178, "Diagnoses", "Primary Diagnoses")
)

$startRow = 150
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $row = $startRow + $i
    $data = $newRows[$i]
    $ws.Cells.Item($row, 1).Value = $data[0]
    $ws.Cells.Item($row, 2).Value = $data[1]
    $ws.Cells.Item($row, 3).Value = $data[2]
    $ws.Cells.Item($row, 4).Value = $data[3]
    $ws.Cells.Item($row, 5).Value = $data[4]
    $ws.Cells.Item($row, 6).Value = $data[5]
    $ws.Cells.Item($row, 7).Value = $data[6]
}
